$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new header row before row 27 (shift old row27..36 down by one)
$ws.Rows.Item(27).Insert()

# New header row 27 (same layout/strings as row 1 header)
$ws.Range("A27").Value = "(Weimar, 10 Clusters)"
$ws.Range("B27").Value = "Floyd-Warshall"
$ws.Range("C27").Value = "All Pairs Dijkstra"
$ws.Range("D27").Value = "All Pairs Dijkstra, parallel 4 cores"
# Row-insert copied the number-format style down onto B27:D27; the header
# cells should carry no explicit style (matches row 1's header formatting).
$ws.Range("B27:D27").ClearFormats()

# Fast UPGMA timing data (row 28 has the label, rows 28-37 have data)
$ws.Range("B28").Value = 53197
$ws.Range("C28").Value = 2409
$ws.Range("D28").Value = 973

$ws.Range("B29").Value = 53321
$ws.Range("C29").Value = 2412
$ws.Range("D29").Value = 975

$ws.Range("B30").Value = 53507
$ws.Range("C30").Value = 2376
$ws.Range("D30").Value = 991

$ws.Range("B31").Value = 53110
$ws.Range("C31").Value = 2394
$ws.Range("D31").Value = 968

$ws.Range("B32").Value = 53411
$ws.Range("C32").Value = 2341
$ws.Range("D32").Value = 973

$ws.Range("B33").Value = 53292
$ws.Range("C33").Value = 2419
$ws.Range("D33").Value = 983

$ws.Range("B34").Value = 53240
$ws.Range("C34").Value = 2442
$ws.Range("D34").Value = 959

$ws.Range("B35").Value = 53248
$ws.Range("C35").Value = 2477
$ws.Range("D35").Value = 969

$ws.Range("B36").Value = 53405
$ws.Range("C36").Value = 2469
$ws.Range("D36").Value = 974

$ws.Range("B37").Value = 53462
$ws.Range("C37").Value = 2386
$ws.Range("D37").Value = 949

# Adjust column A width (no longer auto bestFit) - target stored width is
# 20.42578125; ColumnWidth=19.6 is the closest input this host's pixel
# rounding will serialize toward that value.
$ws.Columns.Item(1).ColumnWidth = 19.6

# Update view: scroll so row 7 is top-left, select D35
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("D35").Select()
